# Refresh the cryptos table on Sheet1 (GitHub Actions data pull):
# update the "Price" (D) and "Volume(1h)" (E) columns for rows 2-51.
#
# Both columns hold plain text in the workbook (prices keep the raw
# scraped "thousand.dot" formatting, e.g. "69.112.98", and percentages
# carry padding spaces), so for any Price value that Excel's automatic
# type-inference would otherwise snap to a Number (dropping significant
# trailing zeros such as "1.00" -> 1, or re-writing "0.0000250" in
# scientific-ish form), the cell is pre-formatted as Text before the
# write and then restored to the default "Normal" style afterwards so
# no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.112.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.814.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "632.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.811.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.453"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.455.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.835.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.108.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.708"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000152"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.15%  "
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.965.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("E38").Value = "  +8.26%  "
$ws.Range("E39").Value = "  +6.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.90%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "157.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("E45").Value = "  +6.54%  "
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  +3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000283"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.65%  "
